$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update USERID (G) and PASSWORD (H) columns for rows 2-5: "Putri"/"bni1234/" -> 32382/"bni1234"
$ws.Range("G2").Value = 32382
$ws.Range("H2").Value = "bni1234"
$ws.Range("F2").Value = "Username : 32382;`nPassword : bni1234;`nKode Group Approval : 9;`nNama Group Approval : Editor"

$ws.Range("G3").Value = 32382
$ws.Range("H3").Value = "bni1234"
$ws.Range("F3").Value = "Username : 32382;`nPassword : bni1234;`nKode Group Approval : 9"

$ws.Range("G4").Value = 32382
$ws.Range("H4").Value = "bni1234"
$ws.Range("F4").Value = "Username : 32382;`nPassword : bni1234;`nKode Group Approval : 9;`nNama Group Approval : Regulator"

$ws.Range("G5").Value = 32382
$ws.Range("H5").Value = "bni1234"
$ws.Range("F5").Value = "Username : 32382;`nPassword : bni1234;`nKode Group Approval : 9"

# Update the active selection to F2 (was D2)
$ws.Range("F2").Select()
